# Fruta / hortaliza, semanal
# Weekly data refresh: insert two new observation rows (new row 5 and row 6)
# at the top of the data block, pushing the previously-existing rows 5-36
# down to rows 7-38.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right above the current row 5 (the first data row),
# shifting all existing data rows down by two.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

function Set-DataRow($r, $vals) {
    $ws.Range("A$r").Value = $vals[0]
    $ws.Range("B$r").Value = $vals[1]
    $ws.Range("C$r").Value = $vals[2]
    $ws.Range("D$r").Value = $vals[3]
    $ws.Range("E$r").Value = $vals[4]
    $ws.Range("F$r").Value = $vals[5]
    $ws.Range("G$r").Value = $vals[6]
    $ws.Range("H$r").Value = $vals[7]
    $ws.Range("I$r").Value = $vals[8]
    $ws.Range("J$r").Value = $vals[9]
    $ws.Range("K$r").Value = $vals[10]
    $ws.Range("L$r").Value = $vals[11]
    $ws.Range("M$r").Value = $vals[12]
    $ws.Range("N$r").Value = $vals[13]
    $ws.Range("O$r").Value = $vals[14]
    $ws.Range("P$r").Value = $vals[15]
    $ws.Range("Q$r").Value = $vals[16]
    $ws.Range("R$r").Value = $vals[17]
    $ws.Range("S$r").Value = $vals[18]
    $ws.Range("T$r").Value = $vals[19]
}

# New row 5: Níspero, Primera, 80 units @ $30000, bandeja de 10 kilos.
Set-DataRow 5 @(
    10, "Vega Modelo de Temuco", "La Araucanía", 45257, 9, "Fruta",
    100104, "Frutos de pepita", 100104004, "Níspero", "Californiana(o)",
    "Primera", 80, 30000, 30000, 30000, "$/bandeja 10 kilos",
    "Provincia de Quillota", 3000, 10
)

# New row 6: Níspero, Primera, 200 units @ $20000, bandeja de 5 kilos.
Set-DataRow 6 @(
    10, "Vega Modelo de Temuco", "La Araucanía", 45257, 9, "Fruta",
    100104, "Frutos de pepita", 100104004, "Níspero", "Californiana(o)",
    "Primera", 200, 20000, 20000, 20000, "$/bandeja 5 kilos",
    "Provincia de Quillota", 4000, 5
)
